# "Added data from Feb 9" -- insert a new observation row (2024-02-09) for the
# ForestHill2023IrrigationFull simulation right after the last existing row
# for that simulation (old row 17 / new row 18), and append a matching new
# observation row for ForestHill2023IrrigationPartial at the end of the
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at 18 (pushes old rows 18..33 down to 19..34).
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the newly inserted row 18 with the Feb-9 ForestHillFull data.
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "ForestHill2023IrrigationFull"
$ws.Range("B18").Value = 45331
$ws.Range("B18").NumberFormat = "d-mmm-yy"
$ws.Range("F18").Value = 5.8926489999999996
$ws.Range("G18").Value = 0.16848393942648088
$ws.Range("H18").NumberFormat = "0.00"
$ws.Range("H18").ClearContents()
$ws.Range("J18").Value = 175.71666666666667
$ws.Range("J18").NumberFormat = "0.00"
$ws.Range("K18").Value = 124.45
$ws.Range("K18").NumberFormat = "0.00"

# The old A18 carried a leftover one-off style (a bold-ish "applyFont" xf)
# that doesn't belong anywhere in the sheet any more -- it shifted down to
# A19 along with the rest of that row's data, so strip it back to the
# default (General) style.
$ws.Range("A19").ClearFormats()

# ---------------------------------------------------------------------------
# 3. Append a new row (35) with the Feb-9 ForestHillPartial data, right
#    after the existing last row (34, formerly 33).
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = "ForestHill2023IrrigationPartial"
$ws.Range("B35").Value = 45331
$ws.Range("B35").NumberFormat = "d-mmm-yy"
$ws.Range("F35").Value = 5.5032705000000002
$ws.Range("G35").Value = 0.016615774246742543
$ws.Range("J35").Value = 165.13333333333333
$ws.Range("K35").Value = 133.63333333333333

# ---------------------------------------------------------------------------
# 4. Keep the _FilterDatabase defined name in sync with the now
#    one-row-taller table.
# ---------------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -like "*_FilterDatabase*") {
        $n.RefersTo = "=CottonObserved!`$A`$1:`$EP`$2578"
    }
}

# ---------------------------------------------------------------------------
# 5. Leave the selection near the newly entered data, like a user would.
# ---------------------------------------------------------------------------
$ws.Range("F18").Select()
